$wb = $excel.ActiveWorkbook

$sheetData = @{}

$sheetData["ALC"] = @(
    @{Row=6; H=552.8; I=322.25; J=706.5; K=966.75; L=2119.5; M=-854.75; N=-2343.5},
    @{Row=39; H=253; I=294.72726; J=100; K=884.18178; L=300; M=-588.18178; N=-892},
    @{Row=118; H=1987; I=149.5; J=3824.5; K=448.5; L=11473.5; M=1208.5; N=-14787.5},
    @{Row=135; H=750; I=750; K=6750; M=-4215},
    @{Row=137; H=1317.4615; J=1647.4286; L=4942.2858; N=-10042.2858}
)

$sheetData["ARM"] = @(
    @{Row=2; H=1988.9166; I=1177.8572; J=3124.4; K=1177.8572; L=3124.4; M=-1064.8572; N=-3350.4},
    @{Row=32; H=4125.48; I=3039.8096; K=3039.8096; M=-2752.8096},
    @{Row=45; H=3627.1365; I=1399.4; K=1399.4; M=-1022.4},
    @{Row=61; H=849.3333; I=849.3333; K=849.3333; M=-637.3333},
    @{Row=116; H=1988.9166; I=1177.8572; J=3124.4; K=1177.8572; L=3124.4; M=1116.1428; N=-7712.4},
    @{Row=122; H=2728.375; I=1450; J=4006.75; K=4350; L=12020.25; M=-1900; N=-16920.25},
    @{Row=132; H=2019.6666; I=1975.2858; J=2175; K=5925.857400000001; L=6525; M=-3395.857400000001; N=-11585},
    @{Row=136; H=849.3333; I=849.3333; K=2547.9999; M=2.000100000000202}
)

$sheetData["BSM"] = @(
    @{Row=3; H=1988.9166; I=1177.8572; J=3124.4; K=1177.8572; L=3124.4; M=-1063.8572; N=-3352.4},
    @{Row=99; H=2669.5715; I=2215.9092; K=2215.9092; M=-717.9092000000001}
)

$sheetData["CRP"] = @(
    @{Row=7; H=1948.5; I=997.6923; J=3714.2856; K=997.6923; L=3714.2856; M=-884.6923; N=-3940.2856},
    @{Row=19; H=1017.5; I=676.25; K=676.25; M=-506.25},
    @{Row=24; H=1017.5; I=676.25; K=676.25; M=-506.25},
    @{Row=31; H=2351.3333; I=1821.6; K=1821.6; M=-1526.6},
    @{Row=34; H=2351.3333; I=1821.6; K=1821.6; M=-1619.6},
    @{Row=99; H=5937.4165; I=6138.778; K=6138.778; M=-4640.778},
    @{Row=122; H=2421.3; I=2246.4443; K=6739.3329; M=-4289.3329},
    @{Row=126; H=5937.4165; I=6138.778; K=18416.334; M=-15946.334},
    @{Row=132; H=2345.1292; I=2411.111; J=1899.75; K=7233.333; L=5699.25; M=-4703.333; N=-10759.25},
    @{Row=134; H=2468.375; I=2106.7144; J=5000; K=6320.1432; L=15000; M=-3785.1432; N=-20070}
)

$sheetData["CUL"] = @(
    @{Row=7; H=485.5; I=0; J=485.5; K=0; L=1456.5; M=$null; N=-1680.5},
    @{Row=12; H=285.25; J=285.25; L=855.75; N=-1201.75},
    @{Row=44; H=175.6; I=172; K=516; M=-118},
    @{Row=51; H=4999; I=0; J=4999; K=0; L=14997; M=$null; N=-15917},
    @{Row=68; H=1668.3334; I=0; J=1668.3334; K=0; L=5005.0002; M=$null; N=-6627.0002},
    @{Row=71; H=1668.3334; I=0; J=1668.3334; K=0; L=15015.0006; M=$null; N=-23127.0006},
    @{Row=80; H=6996.9165; I=5393.2; J=8142.4287; K=16179.6; L=24427.2861; M=-15243.6; N=-26299.2861},
    @{Row=83; H=6996.9165; I=5393.2; J=8142.4287; K=48538.8; L=73281.85830000001; M=-43858.8; N=-82641.85830000001},
    @{Row=137; H=3001.8333; I=990; K=2970; M=2130}
)

$sheetData["GSM"] = @(
    @{Row=3; H=5000000; I=5000000; J=0; K=5000000; L=0; M=-4999884; N=$null},
    @{Row=46; H=16039.5; J=30000; L=30000; N=-30312},
    @{Row=80; H=8301.5; I=4597.5; J=12005.5; K=4597.5; L=12005.5; M=-3599.5; N=-14001.5},
    @{Row=83; H=8301.5; I=4597.5; J=12005.5; K=22987.5; L=60027.5; M=-17995.5; N=-70011.5},
    @{Row=102; H=963.4666999999999; I=788.5; J=1313.4; K=788.5; L=1313.4; M=833.5; N=-4557.4},
    @{Row=107; H=3423.5; I=567.9; J=6993; K=567.9; L=6993; M=1352.1; N=-10833},
    @{Row=122; H=4442.8887; J=3504; L=10512; N=-15412},
    @{Row=126; H=4097.8335; J=4147.25; L=12441.75; N=-17381.75},
    @{Row=132; H=2540.9; I=2540.9; K=7622.700000000001; M=-5092.700000000001}
)

$sheetData["LTW"] = @(
    @{Row=40; H=3000.8; I=3000.8; K=3000.8; M=-2864.8},
    @{Row=100; H=2949.6; I=2949.6; K=2949.6; M=-2408.6},
    @{Row=136; H=4833.3335; I=4833.3335; K=14500.0005; M=-11950.0005}
)

$sheetData["WVR"] = @(
    @{Row=122; H=3300; I=3300; K=9900; M=-7450},
    @{Row=132; H=1665.1428; I=1665.1428; K=4995.428400000001; M=-2465.428400000001}
)

$colOrder = @("H","I","J","K","L","M","N")

$totalCells = 0
foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowData in $sheetData[$sheetName]) {
        $r = $rowData.Row
        foreach ($col in $colOrder) {
            if ($rowData.ContainsKey($col)) {
                $addr = "$col$r"
                $val = $rowData[$col]
                if ($val -eq $null) {
                    $ws.Range($addr).ClearContents()
                } else {
                    $ws.Range($addr).Value = $val
                }
                $totalCells = $totalCells + 1
            }
        }
    }
}

Write-Host "Updated cells:" $totalCells
